$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.361.62"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "'2.574.24"
$ws.Range("E3").Value = "  -2.87%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'542.49"
$ws.Range("D6").Value = "'142.95"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "'0.0999"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "'3.024.31"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "'58.270.62"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "'20.54"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").Value = "'2.565.41"
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "'4.46"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "'333.75"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("D20").Value = "'9.99"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").Value = "'6.10"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'66.34"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "'0.421"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "'0.158"
$ws.Range("E26").Value = "  -5.37%  "
$ws.Range("D27").Value = "'7.04"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'0.0₃0731"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "'5.95"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").Value = "'153.53"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "'18.91"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'3.89"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").Value = "'0.849"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("D36").Value = "'1.09"
$ws.Range("E36").Value = "  -5.05%  "
$ws.Range("D37").Value = "'0.816"
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("D39").Value = "'3.57"
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "'278.68"
$ws.Range("E40").Value = "  -5.85%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.587"
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'10.62"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "'0.0941"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "'0.0529"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").Value = "'18.40"
$ws.Range("E46").Value = "  -5.15%  "
$ws.Range("D47").Value = "'0.0226"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "'1.901.14"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("D49").Value = "'17.82"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").Value = "'4.38"
$ws.Range("E50").Value = "  -3.73%  "
$ws.Range("D51").Value = "'108.70"
$ws.Range("E51").Value = "  -1.92%  "
